$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1476.8334
$ws.Range("I8").Value = 1724.2
$ws.Range("J8").Value = 240
$ws.Range("K8").Value = 5172.6
$ws.Range("L8").Value = 720
$ws.Range("M8").Value = -5033.6
$ws.Range("N8").Value = -998
$ws.Range("H51").Value = 119832.5
$ws.Range("J51").Value = 119832.5
$ws.Range("L51").Value = 119832.5
$ws.Range("N51").Value = -120800.5
$ws.Range("H58").Value = 3282.5557
$ws.Range("I58").Value = 568
$ws.Range("K58").Value = 1704
$ws.Range("M58").Value = -1554
$ws.Range("H132").Value = 966.6786
$ws.Range("I132").Value = 966.6786
$ws.Range("K132").Value = 2900.0358
$ws.Range("M132").Value = -370.0357999999997
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 334.75
$ws.Range("I5").Value = 334.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 334.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -222.75
$ws.Range("N5").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H55").Value = 10048
$ws.Range("I55").Value = 10048
$ws.Range("K55").Value = 10048
$ws.Range("M55").Value = -9733
$ws.Range("H107").Value = 44999.5
$ws.Range("J107").Value = 44999.5
$ws.Range("L107").Value = 44999.5
$ws.Range("N107").Value = -52679.5
$ws.Range("H124").Value = 55484.2
$ws.Range("J124").Value = 55484.2
$ws.Range("L124").Value = 55484.2
$ws.Range("N124").Value = -65304.2
$ws.Range("H125").Value = 70714.5
$ws.Range("J125").Value = 70714.5
$ws.Range("L125").Value = 70714.5
$ws.Range("N125").Value = -80554.5
$ws.Range("H137").Value = 291000
$ws.Range("J137").Value = 291000
$ws.Range("L137").Value = 291000
$ws.Range("N137").Value = -301200

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 334.75
$ws.Range("I4").Value = 334.75
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 334.75
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -219.75
$ws.Range("N4").ClearContents()
$ws.Range("H25").Value = 287.5
$ws.Range("I25").Value = 287.5
$ws.Range("K25").Value = 287.5
$ws.Range("M25").Value = -52.5
$ws.Range("H86").Value = 2208.9285
$ws.Range("I86").Value = 2266.2727
$ws.Range("K86").Value = 2266.2727
$ws.Range("M86").Value = -1143.2727
$ws.Range("H89").Value = 2208.9285
$ws.Range("I89").Value = 2266.2727
$ws.Range("K89").Value = 11331.3635
$ws.Range("M89").Value = -5715.363499999999
$ws.Range("H107").Value = 2278.8572
$ws.Range("I107").Value = 2242
$ws.Range("K107").Value = 2242
$ws.Range("M107").Value = -322
$ws.Range("H134").Value = 57985.25
$ws.Range("I134").Value = 6449.5
$ws.Range("K134").Value = 19348.5
$ws.Range("M134").Value = -16813.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1989.6364
$ws.Range("I7").Value = 227.57143
$ws.Range("J7").Value = 5073.25
$ws.Range("K7").Value = 227.57143
$ws.Range("L7").Value = 5073.25
$ws.Range("M7").Value = -114.57143
$ws.Range("N7").Value = -5299.25
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H25").Value = 7743.2354
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 7743.2354
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 7743.2354
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -8091.2354
$ws.Range("H50").Value = 75000
$ws.Range("J50").Value = 75000
$ws.Range("L50").Value = 75000
$ws.Range("N50").Value = -76250
$ws.Range("H51").Value = 37499.875
$ws.Range("I51").Value = 18199.8
$ws.Range("K51").Value = 18199.8
$ws.Range("M51").Value = -17463.8
$ws.Range("H60").Value = 72416
$ws.Range("J60").Value = 69999
$ws.Range("L60").Value = 69999
$ws.Range("N60").Value = -71021
$ws.Range("H61").Value = 37499.875
$ws.Range("I61").Value = 18199.8
$ws.Range("K61").Value = 18199.8
$ws.Range("M61").Value = -17851.8
$ws.Range("H62").Value = 3075
$ws.Range("J62").Value = 3750
$ws.Range("L62").Value = 3750
$ws.Range("N62").Value = -4998
$ws.Range("H65").Value = 3075
$ws.Range("J65").Value = 3750
$ws.Range("L65").Value = 18750
$ws.Range("N65").Value = -24990
$ws.Range("H74").Value = 80000
$ws.Range("J74").Value = 80000
$ws.Range("L74").Value = 80000
$ws.Range("N74").Value = -81748
$ws.Range("H77").Value = 80000
$ws.Range("J77").Value = 80000
$ws.Range("L77").Value = 240000
$ws.Range("N77").Value = -248736
$ws.Range("H86").Value = 4890.6
$ws.Range("I86").Value = 4989.5
$ws.Range("J86").Value = 4824.6665
$ws.Range("K86").Value = 4989.5
$ws.Range("L86").Value = 4824.6665
$ws.Range("M86").Value = -3866.5
$ws.Range("N86").Value = -7070.6665
$ws.Range("H89").Value = 4890.6
$ws.Range("I89").Value = 4989.5
$ws.Range("J89").Value = 4824.6665
$ws.Range("K89").Value = 24947.5
$ws.Range("L89").Value = 24123.3325
$ws.Range("M89").Value = -19331.5
$ws.Range("N89").Value = -35355.3325
$ws.Range("H132").Value = 2666.5
$ws.Range("I132").Value = 2666.5
$ws.Range("K132").Value = 7999.5
$ws.Range("M132").Value = -5469.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 85000
$ws.Range("J37").Value = 85000
$ws.Range("L37").Value = 255000
$ws.Range("N37").Value = -255224
$ws.Range("H129").Value = 37147840
$ws.Range("J129").Value = 55721308
$ws.Range("L129").Value = 167163924
$ws.Range("N129").Value = -167173924
$ws.Range("H140").Value = 1252.4375
$ws.Range("I140").Value = 1252.4375
$ws.Range("K140").Value = 3757.3125
$ws.Range("M140").Value = 1422.6875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 913.1923
$ws.Range("I97").Value = 768.4583
$ws.Range("K97").Value = 768.4583
$ws.Range("M97").Value = -272.4583
$ws.Range("H132").Value = 62503948
$ws.Range("I132").Value = 76926824
$ws.Range("K132").Value = 230780472
$ws.Range("M132").Value = -230777942

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1823
$ws.Range("I22").Value = 2082.75
$ws.Range("J22").Value = 1433.375
$ws.Range("K22").Value = 2082.75
$ws.Range("L22").Value = 1433.375
$ws.Range("M22").Value = -1787.75
$ws.Range("N22").Value = -2023.375
$ws.Range("H27").Value = 1823
$ws.Range("I27").Value = 2082.75
$ws.Range("J27").Value = 1433.375
$ws.Range("K27").Value = 2082.75
$ws.Range("L27").Value = 1433.375
$ws.Range("M27").Value = -1975.75
$ws.Range("N27").Value = -1647.375
$ws.Range("H132").Value = 107988.63
$ws.Range("I132").Value = 86051.25
$ws.Range("J132").Value = 145595.58
$ws.Range("K132").Value = 258153.75
$ws.Range("L132").Value = 436786.74
$ws.Range("M132").Value = -255623.75
$ws.Range("N132").Value = -441846.74

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 33980
$ws.Range("J109").Value = 33980
$ws.Range("L109").Value = 33980
$ws.Range("N109").Value = -36754
